$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 115.36364
$ws.Range("I9").Value = 96.42856999999999
$ws.Range("K9").Value = 96.42856999999999
$ws.Range("M9").Value = 72.57143000000001

$ws.Range("H62").Value = 116673390
$ws.Range("I62").Value = 50009700
$ws.Range("K62").Value = 50009700
$ws.Range("M62").Value = -50009076

$ws.Range("H65").Value = 116673390
$ws.Range("I65").Value = 50009700
$ws.Range("K65").Value = 250048500
$ws.Range("M65").Value = -250045380

$ws.Range("H98").Value = 22238130
$ws.Range("I98").Value = 8001595.5
$ws.Range("K98").Value = 8001595.5
$ws.Range("M98").Value = -8000097.5

$ws.Range("H122").Value = 22238130
$ws.Range("I122").Value = 8001595.5
$ws.Range("K122").Value = 24004786.5
$ws.Range("M122").Value = -24002336.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 8337672
$ws.Range("I32").Value = 4211.604
$ws.Range("K32").Value = 4211.604
$ws.Range("M32").Value = -3924.604

$ws.Range("H122").Value = 4125.5386
$ws.Range("I122").Value = 3636
$ws.Range("K122").Value = 10908
$ws.Range("M122").Value = -8458

$ws.Range("H132").Value = 34266572
$ws.Range("I132").Value = 36122070
$ws.Range("K132").Value = 108366210
$ws.Range("M132").Value = -108363680

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1902.94
$ws.Range("I86").Value = 1939.3298
$ws.Range("J86").Value = 1332.8334
$ws.Range("K86").Value = 1939.3298
$ws.Range("L86").Value = 1332.8334
$ws.Range("M86").Value = -816.3298
$ws.Range("N86").Value = -3578.8334

$ws.Range("H89").Value = 1902.94
$ws.Range("I89").Value = 1939.3298
$ws.Range("J89").Value = 1332.8334
$ws.Range("K89").Value = 9696.648999999999
$ws.Range("L89").Value = 6664.166999999999
$ws.Range("M89").Value = -4080.648999999999
$ws.Range("N89").Value = -17896.167

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 304.89474
$ws.Range("I22").Value = 246.38461
$ws.Range("J22").Value = 431.66666
$ws.Range("K22").Value = 246.38461
$ws.Range("L22").Value = 431.66666
$ws.Range("M22").Value = 103.61539
$ws.Range("N22").Value = -1131.66666

$ws.Range("H54").Value = 19500
$ws.Range("J54").Value = 19500
$ws.Range("L54").Value = 19500
$ws.Range("N54").Value = -20816

$ws.Range("H58").Value = 1280417.4
$ws.Range("I58").Value = 1731020.2
$ws.Range("J58").Value = 3709.3333
$ws.Range("K58").Value = 1731020.2
$ws.Range("L58").Value = 3709.3333
$ws.Range("M58").Value = -1730817.2
$ws.Range("N58").Value = -4115.3333

$ws.Range("H122").Value = 10590.444
$ws.Range("I122").Value = 11601.75
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 34805.25
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -32355.25
$ws.Range("N122").Value = -12400

$ws.Range("H132").Value = 2275095.8
$ws.Range("I132").Value = 5556755
$ws.Range("J132").Value = 3177.8462
$ws.Range("K132").Value = 16670265
$ws.Range("L132").Value = 9533.5386
$ws.Range("M132").Value = -16667735
$ws.Range("N132").Value = -14593.5386

$ws.Range("H134").Value = 2002791.9
$ws.Range("I134").Value = 1989
$ws.Range("J134").Value = 10006004
$ws.Range("K134").Value = 5967
$ws.Range("L134").Value = 30018012
$ws.Range("M134").Value = -3432
$ws.Range("N134").Value = -30023082

$ws.Range("H136").Value = 1280417.4
$ws.Range("I136").Value = 1731020.2
$ws.Range("J136").Value = 3709.3333
$ws.Range("K136").Value = 5193060.6
$ws.Range("L136").Value = 11127.9999
$ws.Range("M136").Value = -5190510.6
$ws.Range("N136").Value = -16227.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 11809.579
$ws.Range("I80").Value = 4000.6667
$ws.Range("J80").Value = 13273.75
$ws.Range("K80").Value = 12002.0001
$ws.Range("L80").Value = 39821.25
$ws.Range("M80").Value = -11066.0001
$ws.Range("N80").Value = -41693.25

$ws.Range("H83").Value = 11809.579
$ws.Range("I83").Value = 4000.6667
$ws.Range("J83").Value = 13273.75
$ws.Range("K83").Value = 36006.0003
$ws.Range("L83").Value = 119463.75
$ws.Range("M83").Value = -31326.0003
$ws.Range("N83").Value = -128823.75

$ws.Range("H131").Value = 13747.718
$ws.Range("J131").Value = 14854.473
$ws.Range("L131").Value = 44563.419
$ws.Range("N131").Value = -54643.419

$ws.Range("H132").Value = 2237.7222
$ws.Range("I132").Value = 1197.3334
$ws.Range("J132").Value = 3278.111
$ws.Range("K132").Value = 10776.0006
$ws.Range("L132").Value = 29502.999
$ws.Range("M132").Value = -8246.000599999999
$ws.Range("N132").Value = -34562.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 106
$ws.Range("I2").Value = 36.5
$ws.Range("J2").Value = 175.5
$ws.Range("K2").Value = 36.5
$ws.Range("L2").Value = 175.5
$ws.Range("M2").Value = 76.5
$ws.Range("N2").Value = -401.5

$ws.Range("H3").Value = 5000008
$ws.Range("I3").Value = 5000008
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 5000008
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -4999892
$ws.Range("N3").ClearContents()

$ws.Range("H7").Value = 4667333.5
$ws.Range("J7").Value = 4667333.5
$ws.Range("L7").Value = 4667333.5
$ws.Range("N7").Value = -4667557.5

$ws.Range("H8").Value = 4667333.5
$ws.Range("J8").Value = 4667333.5
$ws.Range("L8").Value = 4667333.5
$ws.Range("N8").Value = -4667611.5

$ws.Range("H9").Value = 3166.6667
$ws.Range("I9").Value = 4500
$ws.Range("J9").Value = 2500
$ws.Range("K9").Value = 4500
$ws.Range("L9").Value = 2500
$ws.Range("M9").Value = -4330
$ws.Range("N9").Value = -2840

$ws.Range("H14").Value = 11688333
$ws.Range("J14").Value = 35000
$ws.Range("L14").Value = 35000
$ws.Range("N14").Value = -35336

$ws.Range("H49").Value = 15800
$ws.Range("J49").Value = 15800
$ws.Range("L49").Value = 15800
$ws.Range("N49").Value = -16168

$ws.Range("H80").Value = 17281.45
$ws.Range("I80").Value = 12463.846
$ws.Range("J80").Value = 26228.428
$ws.Range("K80").Value = 12463.846
$ws.Range("L80").Value = 26228.428
$ws.Range("M80").Value = -11465.846
$ws.Range("N80").Value = -28224.428

$ws.Range("H83").Value = 17281.45
$ws.Range("I83").Value = 12463.846
$ws.Range("J83").Value = 26228.428
$ws.Range("K83").Value = 62319.23
$ws.Range("L83").Value = 131142.14
$ws.Range("M83").Value = -57327.23
$ws.Range("N83").Value = -141126.14

$ws.Range("H122").Value = 9806675
$ws.Range("I122").Value = 2483.3333
$ws.Range("J122").Value = 33336734
$ws.Range("K122").Value = 7449.999899999999
$ws.Range("L122").Value = 100010202
$ws.Range("M122").Value = -4999.999899999999
$ws.Range("N122").Value = -100015102

$ws.Range("H132").Value = 25766524
$ws.Range("I132").Value = 23001040
$ws.Range("J132").Value = 29223380
$ws.Range("K132").Value = 69003120
$ws.Range("L132").Value = 87670140
$ws.Range("M132").Value = -69000590
$ws.Range("N132").Value = -87675200

$ws.Range("H141").Value = 31577.6
$ws.Range("J141").Value = 31577.6
$ws.Range("L141").Value = 31577.6
$ws.Range("N141").Value = -41937.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 30409.5
$ws.Range("I93").Value = 18673.715
$ws.Range("J93").Value = 46839.6
$ws.Range("K93").Value = 18673.715
$ws.Range("L93").Value = 46839.6
$ws.Range("M93").Value = -17425.715
$ws.Range("N93").Value = -49335.6

$ws.Range("H136").Value = 33631044
$ws.Range("I136").Value = 11785826
$ws.Range("K136").Value = 35357478
$ws.Range("M136").Value = -35354928

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 14074.8
$ws.Range("I3").Value = 592.5
$ws.Range("J3").Value = 68004
$ws.Range("K3").Value = 592.5
$ws.Range("L3").Value = 68004
$ws.Range("M3").Value = -478.5
$ws.Range("N3").Value = -68232

$ws.Range("H8").Value = 11999
$ws.Range("I8").Value = 11999
$ws.Range("K8").Value = 11999
$ws.Range("M8").Value = -11859

$ws.Range("H9").Value = 25000
$ws.Range("I9").Value = 25000
$ws.Range("K9").Value = 25000
$ws.Range("M9").Value = -24860

$ws.Range("H14").Value = 0
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()

$ws.Range("H122").Value = 1646.9697
$ws.Range("I122").Value = 1555
$ws.Range("J122").Value = 1988.5714
$ws.Range("K122").Value = 4665
$ws.Range("L122").Value = 5965.7142
$ws.Range("M122").Value = -2215
$ws.Range("N122").Value = -10865.7142

$ws.Range("H126").Value = 31250648
$ws.Range("I126").Value = 35714884
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 107144652
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = -107142182
$ws.Range("N126").Value = -7940

$ws.Range("H132").Value = 386949.22
$ws.Range("I132").Value = 578629.8
$ws.Range("J132").Value = 3588
$ws.Range("K132").Value = 1735889.4
$ws.Range("L132").Value = 10764
$ws.Range("M132").Value = -1733359.4
$ws.Range("N132").Value = -15824

$ws.Range("H141").Value = 65000
$ws.Range("J141").Value = 65000
$ws.Range("L141").Value = 65000
$ws.Range("N141").Value = -75360
